# Fix typo in title text on slides 3, 4, 5 and 6:
# "What are these features important?" -> "Why are these features important?"

$p = $ppt.ActivePresentation

$slideIndexes = @(3, 4, 5, 6)

foreach ($idx in $slideIndexes) {
    $s = $p.Slides.Item($idx)
    foreach ($shape in $s.Shapes) {
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq "What are these features important?") {
                $tr.Text = "Why are these features important?"
            }
        }
    }
}
